$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells (AD1:AF1) - "Wins", "Losses", "Ties"
$ws.Cells.Item(1, 30).Value = "Wins"
$ws.Cells.Item(1, 31).Value = "Losses"
$ws.Cells.Item(1, 32).Value = "Ties"

# Copy the header style (bold, centered, bordered) from an existing header cell (AC1)
# onto the new header cells so they match the rest of row 1.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Fill in the team record (Wins/Losses/Ties) for every player row.
for ($r = 2; $r -le 44; $r++) {
    $ws.Cells.Item($r, 30).Value = 81
    $ws.Cells.Item($r, 31).Value = 81
    $ws.Cells.Item($r, 32).Value = 0
}
